$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

$companies = @("company","臺銀人壽","臺銀人壽","富邦人壽","國泰人壽","台灣人壽","國華人壽","富邦人壽","國泰人壽","國泰人壽","國泰人壽","新光人壽","國泰人壽","新光人壽","幸福人壽","中華郵政","台灣人壽")
$names = @("name","美麗人生萬能保險","牛轉錢坤萬能保險","豐年養老保險","卓越理財變額萬能壽險","掌握人生變額保險","定期終身保險","吉百利終身壽險","富貴保本三福終身險","鍾愛一生重大疾病險","雙喜年年終身險","長安終身壽險","鍾愛一生313终身險","威利長福保險","白金人生養老保險","常春增額還本保險","長期看護終身險")
$owners = @("owner","王進士","王進士","周麗容","王進士","王進士","王進士","周麗容","周麗容","周麗容","周麗容","周麗容","周麗容","周麗容","周麗容","周麗容","周麗容")
$propcat = @("property_category","insurance","insurance","insurance","insurance","insurance","insurance","insurance","insurance","insurance","insurance","insurance","insurance","insurance","insurance","insurance","insurance")
$category = @("category","normal","normal","normal","normal","normal","normal","normal","normal","normal","normal","normal","normal","normal","normal","normal","normal")
$date = @("date","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12","2013-12-12")
$legname = @("legislator_name","王進士","王進士","王進士","王進士","王進士","王進士","王進士","王進士","王進士","王進士","王進士","王進士","王進士","王進士","王進士","王進士")
$legid = @("legislator_id",1701,1701,1701,1701,1701,1701,1701,1701,1701,1701,1701,1701,1701,1701,1701,1701)
$srcfile = @("source_file","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21","tmp93a21")
$indexcol = @("index",106,107,108,109,110,111,112,113,115,116,117,118,119,120,121,122)

for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 2).Value = $companies[$i]
}
for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 3).Value = $names[$i]
}
for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 4).Value = $owners[$i]
}
for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 5).Value = $propcat[$i]
}
for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 6).Value = $category[$i]
}
for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 7).Value = $date[$i]
}
for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 8).Value = $legname[$i]
}
for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 9).Value = $legid[$i]
}
for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 10).Value = $srcfile[$i]
}
for ($i = 0; $i -lt 17; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 11).Value = $indexcol[$i]
}
